$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (it sat between "processing" and
#    " them too." in the Project Steps section) -- it will be re-added at
#    the new last-edit location below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Super-project" -> "Super-Project" (capitalize the P), splitting the
#    text into three runs ("Super-", "P", "roject") with a new "_GoBack"
#    bookmark sitting right after the "P" (the simulated edit point).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Super-project", $true, $false, $false, $false, $false, $true, 1, $false, "Super-Project", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("Super-Project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$superStart = $rng2.Start

# Temporary bookmark right after "Super-" forces a run split there; it gets
# removed again once the real "_GoBack" bookmark (after "Super-P") has been
# inserted, leaving three runs: "Super-" | "P" | "roject".
$tempSplitPoint = $superStart + 6
$d.Bookmarks.Add("TempSplitMark", $d.Range($tempSplitPoint, $tempSplitPoint)) | Out-Null

$goBackPoint = $superStart + 7
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPoint, $goBackPoint)) | Out-Null

$d.Bookmarks.Item("TempSplitMark").Delete()

# ---------------------------------------------------------------------------
# 3) Add a "Date & " run before the "Time" heading run.
# ---------------------------------------------------------------------------
$timeRng = $d.Content
$timeRng.Find.Execute("Time", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$timeStart = $timeRng.Start

# Same temp-bookmark trick so the inserted text lands in its own run instead
# of being merged into the following "Time" run.
$d.Bookmarks.Add("TempSplitMark2", $d.Range($timeStart, $timeStart)) | Out-Null
$insertRng = $d.Range($timeStart, $timeStart)
$insertRng.InsertBefore("Date & ")
$d.Bookmarks.Item("TempSplitMark2").Delete()

# ---------------------------------------------------------------------------
# 4) Heading 3 style: bump the run size to 10pt (w:sz 20).
# ---------------------------------------------------------------------------
$heading3 = $d.Styles.Item("Heading3")
$heading3.Font.Size = 10
